$d = $word.ActiveDocument

# Paragraph 2 (subtitle): split author names into separate runs with proofErr spell-check markers
$p = $d.Paragraphs(2).Range
$xml = @'
<w:p w14:paraId="0A7768B3" w14:textId="5C2DFA45" w:rsidR="006557B2" w:rsidRDefault="006557B2" w:rsidP="006557B2"><w:pPr><w:pStyle w:val="Subtitle"/></w:pPr><w:r><w:t xml:space="preserve">18740: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Shravani</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Dhote</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Simrit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Kaur, Vins Sharma</w:t></w:r></w:p>
'@
[void]$p.InsertXML($xml)

# Paragraph 9: split " objects in " -> " obj" + "ects" (proofErr) + " in "
$p = $d.Paragraphs(9).Range
$xml = @'
<w:p w14:paraId="74FB686F" w14:textId="77777777" w:rsidR="00473F3D" w:rsidRDefault="00AF3337" w:rsidP="00AF3337"><w:r><w:t xml:space="preserve">The base code emulates an 8-way memory structure. In order to modify this, we noticed that the </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>cache_lines[]</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve"> list was a bounded list of </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>Line[]</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve"> lists, which themselves were bounded to </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>8</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve"> obj</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ects</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>need_eviction[]</m:t></m:r></m:oMath><w:r w:rsidR="00473F3D"><w:t>.</w:t></w:r></w:p>
'@
[void]$p.InsertXML($xml)

# Paragraph 10: split "etc" and "sizings" out into their own proofErr-wrapped runs
$p = $d.Paragraphs(10).Range
$xml = @'
<w:p w14:paraId="115D6C0B" w14:textId="777DFC9F" w:rsidR="00BD5C91" w:rsidRDefault="00473F3D" w:rsidP="00AF3337"><w:r><w:t xml:space="preserve">Our implementation strategy was simple – In order to add in core-locked ways </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>without</w:t></w:r><w:r><w:t xml:space="preserve"> modifying the cache’s overall functionality, </w:t></w:r><w:r w:rsidR="00BD5C91"><w:t xml:space="preserve">associativity, block sizing, set sizing, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, and still reuse the base code, we cut down the number of ways to </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>2</m:t></m:r></m:oMath><w:r w:rsidR="00BD5C91"><w:t xml:space="preserve"> and duplicated the </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>cache_lines</m:t></m:r></m:oMath><w:r w:rsidR="00BD5C91"><w:t xml:space="preserve"> list of cache line sets by </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>4</m:t></m:r></m:oMath><w:r w:rsidR="00BD5C91"><w:t xml:space="preserve">. This gave us </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>8</m:t></m:r></m:oMath><w:r w:rsidR="00BD5C91"><w:t xml:space="preserve">-way associativity overall, where every set of </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>2</m:t></m:r></m:oMath><w:r w:rsidR="00BD5C91"><w:t xml:space="preserve"> ways was specified to a particular core</w:t></w:r><w:r w:rsidR="00F83A0D"><w:t xml:space="preserve">, without changing the number of sets, block </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sizings</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, or any significant functionality.</w:t></w:r></w:p>
'@
[void]$p.InsertXML($xml)

# Paragraph 15: split "lock up" and "performance" similarly
$p = $d.Paragraphs(15).Range
$xml = @'
<w:p w14:paraId="3AA85D62" w14:textId="77777777" w:rsidR="00DE0690" w:rsidRDefault="00A96175" w:rsidP="004C522E"><w:r><w:t xml:space="preserve">Interestingly, we found that (by functional specification from the PDF) it was possible for all cores to be locked up with no memory access until the next </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>10,000</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve"> cycle quanta.</w:t></w:r><w:r w:rsidR="00DE0690"><w:t xml:space="preserve"> If one core performed </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>5</m:t></m:r></m:oMath><w:r w:rsidR="00DE0690"><w:t xml:space="preserve"> requests back-to-back, it would be blacklisted, and the other </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>3</m:t></m:r></m:oMath><w:r w:rsidR="00DE0690"><w:t xml:space="preserve"> cores were capable of sending requests much faster, increasing their chances of being blacklisted, etc. As such, the “blacklisting” part of BLISS would stop being relevant after three cores l</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ock</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> up – After the next </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>5</m:t></m:r></m:oMath><w:r w:rsidR="00DE0690"><w:t xml:space="preserve"> consecutive requests from one core, all cores would be blacklisted, and would be served in a first-come first-</w:t></w:r><w:r w:rsidR="00DE0690"><w:lastRenderedPageBreak/><w:t xml:space="preserve">serve basis until the start of the next </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>10,000</m:t></m:r></m:oMath><w:r w:rsidR="00DE0690"><w:t xml:space="preserve"> cycle quanta (and the blacklisting would have no significant effect on pe</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rformance</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>).</w:t></w:r></w:p>
'@
[void]$p.InsertXML($xml)

# New trailing paragraph after the final paragraph
$last = $d.Paragraphs($d.Paragraphs.Count).Range
$last.InsertParagraphAfter()
$newPara = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$newPara.Text = "The results of this showed improvements on our prior Equity scheduler overall, "

